$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM rows (Part Name / Copy number / Note) for rows 13-16.
# Values are written in this particular order so new shared-string
# entries land in the same sequence as the source workbook.
$ws.Range("B16").Value = "m3 nut"
$ws.Range("B14").Value = "m3 anti-slip nut"
$ws.Range("D13").Value = "metal"
$ws.Range("B15").Value = "m3 x 8"
$ws.Range("B13").Value = "m3 x 28"

$ws.Range("D14").Value = "metal"
$ws.Range("D15").Value = "metal"
$ws.Range("D16").Value = "metal"

$ws.Range("C13").Value = 12
$ws.Range("C14").Value = 20
$ws.Range("C15").Value = 8
$ws.Range("C16").Value = 4

# Set column widths to match new content widths (Excel bestFit auto-sizing)
$ws.Columns("B").ColumnWidth = 14.4
$ws.Columns("D").ColumnWidth = 14.4

# Selection moves to D20
$ws.Range("D20").Select()
